# Improved World Bank Data Fetch: Checks for Missing Countries and Runs Them
# & Handles Countries That's Missing All Data
#
# Updates the "Corruption Index (2024)" column (D) for the countries whose
# values were re-fetched / corrected, then leaves the selection where the
# author's cursor ended up (D53 - Turkey).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29 - Kazakhstan
$ws.Range("D29").Value = 40
# Row 31 - South Korea
$ws.Range("D31").Value = 64
# Row 35 - Mongolia
$ws.Range("D35").Value = 33
# Row 48 - Russia
$ws.Range("D48").Value = 22
# Row 53 - Turkey
$ws.Range("D53").Value = 34
# Row 54 - Taiwan
$ws.Range("D54").Value = 67

# Match the final cursor position recorded in the workbook view.
$ws.Range("D53").Select() | Out-Null
